$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3)
$ws.Range("B3").Value = "6.0.0"

# Update Date value (row 8)
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (row 9) was empty -> now "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was a duplicated "Contact" / "No display for ContactDetail" row.
# It becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# The old row 11 duplicated "Contact" / "No display for ContactDetail" entirely
# and is no longer needed now that row 10 carries the Jurisdiction data;
# remove it so everything below shifts up by one row.
$ws.Range("A11").EntireRow.Delete()

# Case Sensitive value (now row 14 after the deletion) was empty -> now "true"
$ws.Range("B14").Value = "true"
